$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sonuclar")

$data = @(
    @(93, "login-functionality;login-with-valid-username-and-password", "failed", "20201225_224205", "chrome"),
    @(94, "country-create-edit-delete-functionality;create-a-country", "failed", "20201225_224222", "chrome"),
    @(95, "country-create-edit-delete-functionality;edit-a-country", "passed", "20201225_224249", "chrome"),
    @(96, "country-create-edit-delete-functionality;delete-a-country", "passed", "20201225_224316", "chrome")
)

$startRow = 94
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
